# Append the Hillcrest Foods / Bakery order lines (rows 3-19) below the
# existing header (row 1) and first data row (row 2). Every column in this
# sheet is stored as text (SKU codes, quantities and prices are all plain
# strings, not numbers), so each value is entered with a leading apostrophe
# to force text entry, and the style is reset to "Normal" afterwards so the
# temporary quote-prefix formatting doesn't linger on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("62034", "M&Ms Bulk", "1", "112.55", "112.55")
    ,@("70328", "Corn Meal - Coarse", "2", "26.49", "52.98")
    ,@("25170", "Raisins - Midget", "2", "70.15", "140.30")
    ,@("33100", "Nuts - Almonds (whole raw)", "1", "116.03", "116.03")
    ,@("33069", "Nuts - Almonds (sliced)", "1", "118.44", "118.44")
    ,@("33053", "Nuts-Almonds Blanched (Slivered)", "1", "113.68", "113.68")
    ,@("33422", "Nuts - Walnut Halves & Pieces", "3", "124.04", "372.12")
    ,@("47170", "Seeds - Sesame", "6", "67.24", "403.44")
    ,@("17126", "Flour - Warrior Clear", "1", "16.24", "16.24")
    ,@("18580", "Flour - Victoria Patent (Org)", "6", "40.02", "240.12")
    ,@("14145", "Chocolate - Unsweetened Dark", "1", "183.00", "183.00")
    ,@("28173", "Sugar - 10X (Powdered)", "3", "45.71", "137.13")
    ,@("72020", "Baking Soda", "1", "37.40", "37.40")
    ,@("70039", "Coconut Swt Flake", "1", "62.37", "62.37")
    ,@("48133", "Spice - Cinnamon (Ground)", "3", "44.20", "132.60")
    ,@("36208", "FRZ Fruit - Mango", "2", "53.93", "107.86")
    ,@("36411", "Rhubarb - FRZ (IQF)", "1", "56.43", "56.43")
)

$row = 3
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = "'" + $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = "'" + $r[2]
    $ws.Cells.Item($row, 4).Value = "'" + $r[3]
    $ws.Cells.Item($row, 5).Value = "'" + $r[4]

    $rowRange = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 5))
    $rowRange.Style = "Normal"

    $row++
}
